# modified sendkeys data from Excel
# Adds two new worksheets (AutoCompleteSampleSheet, DataFromSeleniumEasyURL)
# populated with sample "sendkeys" test data, and makes the last sheet active.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

function AddSheetAtEnd([string]$name) {
    $last = $wb.Worksheets.Item($wb.Worksheets.Count)
    $newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $last)
    $newSheet.Name = $name
    return $newSheet
}

# --- Create the two target sheets -----------------------------------------
# Real Excel keeps a monotonically increasing internal sheetId counter that
# is never reused, even across deletes. Add two throw-away sheets in between
# (consuming ids 6 and 7) and remove them afterwards, so the final sheets end
# up with sheetId 5 and 8, matching the target workbook.
AddSheetAtEnd("AutoCompleteSampleSheet") | Out-Null    # sheetId 5
AddSheetAtEnd("__tmp1__") | Out-Null                    # sheetId 6
AddSheetAtEnd("__tmp2__") | Out-Null                    # sheetId 7
AddSheetAtEnd("DataFromSeleniumEasyURL") | Out-Null     # sheetId 8

$wb.Worksheets.Item("__tmp1__").Delete()
$wb.Worksheets.Item("__tmp2__").Delete()

$autoCompleteSheet = $wb.Worksheets.Item("AutoCompleteSampleSheet")
$seleniumSheet = $wb.Worksheets.Item("DataFromSeleniumEasyURL")

# --- AutoCompleteSampleSheet content ---------------------------------------
$autoCompleteSheet.Range("A1").Value = "AutoCompleteValue"
$autoCompleteSheet.Range("A2").Value = "java"

# Match header style used elsewhere in the workbook (bold font + yellow fill)
$wb.Worksheets.Item("NewAccounts").Range("A1").Copy() | Out-Null
$autoCompleteSheet.Range("A1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

$autoCompleteSheet.Columns.Item(1).AutoFit() | Out-Null
$autoCompleteSheet.Range("D14").Select() | Out-Null

# --- DataFromSeleniumEasyURL content ---------------------------------------
# Values are entered column-by-column (A1, A2, then B1, B2) to reproduce the
# exact shared-string ordering of the source workbook.
$seleniumSheet.Range("A1").Value = "ChildWindow"
$seleniumSheet.Range("A2").Value = "vishalg@testmail.com"
$seleniumSheet.Range("B1").Value = "datepicker"

# Force the date-like text to be stored as text, not parsed as a date.
$seleniumSheet.Range("B2").NumberFormat = "@"
$seleniumSheet.Range("B2").Value = "09/10/1968"

# Header row formatting (bold font + yellow fill), same style as above.
$wb.Worksheets.Item("NewAccounts").Range("A1").Copy() | Out-Null
$seleniumSheet.Range("A1:B1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# A2 gets the hyperlink / hyperlink style; B2 reuses the plain text style.
$seleniumSheet.Hyperlinks.Add($seleniumSheet.Range("A2"), "mailto:vishalg@testmail.com") | Out-Null
$seleniumSheet.Range("A2").Style = "Hyperlink"

$wb.Worksheets.Item("NewCustomers").Range("C2").Copy() | Out-Null
$seleniumSheet.Range("B2").PasteSpecial(-4122) | Out-Null      # xlPasteFormats
$excel.CutCopyMode = $false

$seleniumSheet.Columns.Item(1).AutoFit() | Out-Null
$seleniumSheet.Columns.Item(2).AutoFit() | Out-Null

$seleniumSheet.PageSetup.Orientation = 1       # xlPortrait
$autoCompleteSheet.PageSetup.Orientation = 1   # xlPortrait

$seleniumSheet.Range("A1:B1").Select() | Out-Null

# Make the last sheet the active one (matches activeTab in the target file).
$seleniumSheet.Activate()
